$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing values per diff
$ws.Range("B8").Value = 2298566
$ws.Range("B12").Value = 2040684
$ws.Range("B15").Value = 1974824
$ws.Range("B16").Value = 1674458
$ws.Range("B19").Value = 2547799
$ws.Range("B20").Value = 2024343
$ws.Range("C20").Value = 1858823

# Add new row 21 (Provisional Occupancy Forecast)
# Copy the date-formatted style from A20 into A21 before setting the value
$ws.Range("A20").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = 43544

$ws.Range("B21").Value = 2091997
$ws.Range("C21").Value = 1685390
$ws.Range("D21").Formula = "=B21+C21/Hoja2!`$A`$2"

# Grow the Tabla1 table to include the new row
$lo = $ws.ListObjects.Item("Tabla1")
$lo.Resize($ws.Range("A1:D21"))

$wb.Save()
